$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: cell address, new text value, whether the column needs the
# "force text" treatment (so Excel does not silently coerce e.g. "1.001"
# into the number 1.001 and lose the original inline-string formatting).
$updates = @(
    @{ Cell = "D2"; Value = "27.553.46"; ForceText = $True }
    @{ Cell = "E2"; Value = "  -0.07%  "; ForceText = $False }
    @{ Cell = "D3"; Value = "1.754.53"; ForceText = $True }
    @{ Cell = "E3"; Value = "  +0.05%  "; ForceText = $False }
    @{ Cell = "D4"; Value = "1.001"; ForceText = $True }
    @{ Cell = "E4"; Value = "  +0.00%  "; ForceText = $False }
    @{ Cell = "D5"; Value = "324.36"; ForceText = $True }
    @{ Cell = "E5"; Value = "  +0.07%  "; ForceText = $False }
    @{ Cell = "D6"; Value = "1.000"; ForceText = $True }
    @{ Cell = "E6"; Value = "  -0.01%  "; ForceText = $False }
    @{ Cell = "D7"; Value = "0.4567"; ForceText = $True }
    @{ Cell = "E7"; Value = "  +2.60%  "; ForceText = $False }
    @{ Cell = "D8"; Value = "0.3558"; ForceText = $True }
    @{ Cell = "E8"; Value = "  -1.73%  "; ForceText = $False }
    @{ Cell = "D9"; Value = "0.07476"; ForceText = $True }
    @{ Cell = "E9"; Value = "  -0.35%  "; ForceText = $False }
    @{ Cell = "D10"; Value = "41.53"; ForceText = $True }
    @{ Cell = "E10"; Value = "  -1.67%  "; ForceText = $False }
    @{ Cell = "D11"; Value = "1.086"; ForceText = $True }
    @{ Cell = "E11"; Value = "  -1.86%  "; ForceText = $False }
    @{ Cell = "E12"; Value = "  -0.01%  "; ForceText = $False }
    @{ Cell = "D13"; Value = "20.76"; ForceText = $True }
    @{ Cell = "E13"; Value = "  +0.28%  "; ForceText = $False }
    @{ Cell = "D14"; Value = "6.015"; ForceText = $True }
    @{ Cell = "E14"; Value = "  -0.51%  "; ForceText = $False }
    @{ Cell = "D15"; Value = "7.176"; ForceText = $True }
    @{ Cell = "E15"; Value = "  -0.19%  "; ForceText = $False }
    @{ Cell = "D16"; Value = "1.757.78"; ForceText = $True }
    @{ Cell = "E16"; Value = "  +0.37%  "; ForceText = $False }
    @{ Cell = "D17"; Value = "94.19"; ForceText = $True }
    @{ Cell = "E17"; Value = "  +1.24%  "; ForceText = $False }
    @{ Cell = "D18"; Value = "0.00001056"; ForceText = $True }
    @{ Cell = "E18"; Value = "  -0.68%  "; ForceText = $False }
    @{ Cell = "D19"; Value = "0.06409"; ForceText = $True }
    @{ Cell = "E19"; Value = "  -0.19%  "; ForceText = $False }
    @{ Cell = "E20"; Value = "  -0.03%  "; ForceText = $False }
    @{ Cell = "E21"; Value = "  +0.19%  "; ForceText = $False }
    @{ Cell = "D22"; Value = "5.749"; ForceText = $True }
    @{ Cell = "E22"; Value = "  -1.60%  "; ForceText = $False }
    @{ Cell = "D23"; Value = "27.606.61"; ForceText = $True }
    @{ Cell = "D24"; Value = "11.20"; ForceText = $True }
    @{ Cell = "E24"; Value = "  -0.57%  "; ForceText = $False }
    @{ Cell = "D25"; Value = "2.082"; ForceText = $True }
    @{ Cell = "E25"; Value = "  -0.85%  "; ForceText = $False }
    @{ Cell = "D26"; Value = "165.78"; ForceText = $True }
    @{ Cell = "E26"; Value = "  +1.78%  "; ForceText = $False }
    @{ Cell = "E27"; Value = "  -1.36%  "; ForceText = $False }
    @{ Cell = "D28"; Value = "1.957.46"; ForceText = $True }
    @{ Cell = "E28"; Value = "  +0.37%  "; ForceText = $False }
    @{ Cell = "D29"; Value = "2.127"; ForceText = $True }
    @{ Cell = "E29"; Value = "  -0.13%  "; ForceText = $False }
    @{ Cell = "D30"; Value = "125.72"; ForceText = $True }
    @{ Cell = "E30"; Value = "  -0.02%  "; ForceText = $False }
    @{ Cell = "D31"; Value = "1.084"; ForceText = $True }
    @{ Cell = "E31"; Value = "  -0.27%  "; ForceText = $False }
    @{ Cell = "D32"; Value = "0.09231"; ForceText = $True }
    @{ Cell = "E32"; Value = "  +2.36%  "; ForceText = $False }
    @{ Cell = "D33"; Value = "3.655"; ForceText = $True }
    @{ Cell = "E33"; Value = "  +0.46%  "; ForceText = $False }
    @{ Cell = "D34"; Value = "5.528"; ForceText = $True }
    @{ Cell = "E34"; Value = "  -0.40%  "; ForceText = $False }
    @{ Cell = "B35"; Value = "Aptos"; ForceText = $False }
    @{ Cell = "C35"; Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; ForceText = $False }
    @{ Cell = "D35"; Value = "11.74"; ForceText = $True }
    @{ Cell = "E35"; Value = "  -3.07%  "; ForceText = $False }
    @{ Cell = "B36"; Value = "VeChain"; ForceText = $False }
    @{ Cell = "C36"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; ForceText = $False }
    @{ Cell = "D36"; Value = "0.02283"; ForceText = $True }
    @{ Cell = "E36"; Value = "  -1.08%  "; ForceText = $False }
    @{ Cell = "E37"; Value = "  -0.31%  "; ForceText = $False }
    @{ Cell = "D38"; Value = "0.06016"; ForceText = $True }
    @{ Cell = "E38"; Value = "  +0.93%  "; ForceText = $False }
    @{ Cell = "D39"; Value = "0.6298"; ForceText = $True }
    @{ Cell = "E39"; Value = "  -0.98%  "; ForceText = $False }
    @{ Cell = "D40"; Value = "4.925"; ForceText = $True }
    @{ Cell = "E40"; Value = "  -0.53%  "; ForceText = $False }
    @{ Cell = "D41"; Value = "1.181"; ForceText = $True }
    @{ Cell = "E41"; Value = "  -1.20%  "; ForceText = $False }
    @{ Cell = "D42"; Value = "1.386"; ForceText = $True }
    @{ Cell = "E42"; Value = "  +0.05%  "; ForceText = $False }
    @{ Cell = "D43"; Value = "7.798"; ForceText = $True }
    @{ Cell = "E43"; Value = "  -0.11%  "; ForceText = $False }
    @{ Cell = "D44"; Value = "13.19"; ForceText = $True }
    @{ Cell = "E44"; Value = "  -0.21%  "; ForceText = $False }
    @{ Cell = "D45"; Value = "3.716"; ForceText = $True }
    @{ Cell = "E45"; Value = "  +0.06%  "; ForceText = $False }
    @{ Cell = "D46"; Value = "0.5851"; ForceText = $True }
    @{ Cell = "E46"; Value = "  -0.49%  "; ForceText = $False }
    @{ Cell = "D47"; Value = "122.07"; ForceText = $True }
    @{ Cell = "E48"; Value = "  -1.35%  "; ForceText = $False }
    @{ Cell = "D49"; Value = "0.06891"; ForceText = $True }
    @{ Cell = "E49"; Value = "  +0.45%  "; ForceText = $False }
    @{ Cell = "D50"; Value = "1.130"; ForceText = $True }
    @{ Cell = "D51"; Value = "71.98"; ForceText = $True }
    @{ Cell = "E51"; Value = "  -0.63%  "; ForceText = $False }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Temporarily mark the cell as Text so Excel stores the literal
        # characters instead of re-parsing them as a number, then drop
        # back to the workbook default ("Normal") style so no stray
        # per-cell formatting is left behind.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
